# EQORE deck trim: keep only the title slide and the methodology slide,
# delete the remaining five slides, and rewrite the surviving text.

$p = $ppt.ActivePresentation

# Remove slides 3-7 (delete from the end so indices of earlier slides stay valid)
for ($i = $p.Slides.Count; $i -ge 3; $i--) {
    $p.Slides.Item($i).Delete()
}

# Slide 1: title + subtitle
$slide1 = $p.Slides.Item(1)
$slide1.Shapes.Item(1).TextFrame.TextRange.Text = "EQORE Funding Deck"
$slide1.Shapes.Item(2).TextFrame.TextRange.Text = "Auto-generated deck"

# Slide 2: title + clear the bullet body
$slide2 = $p.Slides.Item(2)
$slide2.Shapes.Item(1).TextFrame.TextRange.Text = "Top Opportunities"
$slide2.Shapes.Item(2).TextFrame.TextRange.Text = ""
